# Applies the cryptos.xlsx data refresh described in the commit:
# "Updated cryptos list on Wed Dec 20 18:27:50 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "44.036.47"
$ws.Range("E2").Value = "  +4.53%  "

# Row 3
$ws.Range("D3").Value = "2.229.31"
$ws.Range("E3").Value = "  +3.32%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.64"
$ws.Range("E5").Value = "  +3.03%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "83.10"
$ws.Range("E6").Value = "  +13.91%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  +1.97%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.610"
$ws.Range("E9").Value = "  +5.43%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.18"
$ws.Range("E10").Value = "  +11.62%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0930"
$ws.Range("E11").Value = "  +2.41%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.07"
$ws.Range("E12").Value = "  +5.31%  "

# Row 13
$ws.Range("E13").Value = "  +2.94%  "

# Row 14
$ws.Range("D14").Value = "2.569.51"
$ws.Range("E14").Value = "  +3.73%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.71"
$ws.Range("E15").Value = "  +4.22%  "

# Row 16
$ws.Range("D16").Value = "2.239.36"
$ws.Range("E16").Value = "  +4.22%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.785"
$ws.Range("E17").Value = "  +3.08%  "

# Row 18
$ws.Range("D18").Value = "43.912.98"
$ws.Range("E18").Value = "  +4.55%  "

# Row 19
$ws.Range("E19").Value = "  +2.65%  "

# Row 20
$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.47"
$ws.Range("E20").Value = "  +1.59%  "

# Row 21
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.04"
$ws.Range("E21").Value = "  +3.83%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.37"
$ws.Range("E22").Value = "  +10.65%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.37"
$ws.Range("E23").Value = "  +3.45%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.26"
$ws.Range("E24").Value = "  -2.36%  "

# Row 25
$ws.Range("E25").Value = "  +0.06%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.78"
$ws.Range("E26").Value = "  +3.65%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "40.61"
$ws.Range("E27").Value = "  +10.84%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.36"
$ws.Range("E28").Value = "  +1.62%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  +2.43%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.23"
$ws.Range("E30").Value = "  +0.35%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.40"
$ws.Range("E31").Value = "  +3.25%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0898"
$ws.Range("E32").Value = "  +12.04%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.68"
$ws.Range("E33").Value = "  +4.04%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.34"
$ws.Range("E34").Value = "  +4.87%  "

# Row 35
$ws.Range("E35").Value = "  +9.19%  "

# Row 36
$ws.Range("E36").Value = "  +2.71%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0367"
$ws.Range("E37").Value = "  +12.10%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.53"
$ws.Range("E38").Value = "  +7.53%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.41"
$ws.Range("E39").Value = "  +14.12%  "

# Row 40
$ws.Range("E40").Value = "  +26.10%  "

# Row 41
$ws.Range("E41").Value = "  +4.64%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "64.02"
$ws.Range("E42").Value = "  +9.67%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.53"
$ws.Range("E43").Value = "  +8.44%  "

# Row 44
$ws.Range("E44").Value = "  +4.22%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.75"
$ws.Range("E45").Value = "  +2.08%  "

# Row 46
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0988"
$ws.Range("E46").Value = "  +3.04%  "

# Row 47
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.39"
$ws.Range("E47").Value = "  +2.56%  "

# Row 48
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.58"
$ws.Range("E48").Value = "  +29.78%  "

# Row 49
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.13"
$ws.Range("E49").Value = "  +4.75%  "

# Row 50
$ws.Range("E50").Value = "  -1.85%  "

# Row 51
$ws.Range("E51").Value = "  +3.93%  "
